$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update D1 text value (shared string changes from "no. telp" to "no.telp")
$ws.Range("D1").Value = "no.telp"

# Apply a Text number format to column D (numFmtId 49 = "@")
$ws.Columns("D").NumberFormat = "@"

# Update selection to G4
$ws.Range("G4").Select()
